$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-15 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-22 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("261÷5=52, 1", $true, $false, $false, $false, $false, $true, 1, $false, "732÷6=122, 0", 2) | Out-Null
$d.Content.Find.Execute("770÷5=154, 0", $true, $false, $false, $false, $false, $true, 1, $false, "925÷8=115, 5", 2) | Out-Null
$d.Content.Find.Execute("613÷6=102, 1", $true, $false, $false, $false, $false, $true, 1, $false, "555÷9=61, 6", 2) | Out-Null
$d.Content.Find.Execute("137÷7=19, 4", $true, $false, $false, $false, $false, $true, 1, $false, "634÷4=158, 2", 2) | Out-Null
$d.Content.Find.Execute("451÷7=64, 3", $true, $false, $false, $false, $false, $true, 1, $false, "380÷8=47, 4", 2) | Out-Null
$d.Content.Find.Execute("805÷2=402, 1", $true, $false, $false, $false, $false, $true, 1, $false, "852÷6=142, 0", 2) | Out-Null
$d.Content.Find.Execute("187÷8=23, 3", $true, $false, $false, $false, $false, $true, 1, $false, "807÷4=201, 3", 2) | Out-Null
$d.Content.Find.Execute("393÷7=56, 1", $true, $false, $false, $false, $false, $true, 1, $false, "178÷7=25, 3", 2) | Out-Null
$d.Content.Find.Execute("764÷3=254, 2", $true, $false, $false, $false, $false, $true, 1, $false, "612÷8=76, 4", 2) | Out-Null
$d.Content.Find.Execute("831÷7=118, 5", $true, $false, $false, $false, $false, $true, 1, $false, "998÷9=110, 8", 2) | Out-Null
$d.Content.Find.Execute("507÷4=126, 3", $true, $false, $false, $false, $false, $true, 1, $false, "568÷8=71, 0", 2) | Out-Null
$d.Content.Find.Execute("805÷9=89, 4", $true, $false, $false, $false, $false, $true, 1, $false, "428÷3=142, 2", 2) | Out-Null
$d.Content.Find.Execute("913÷3=304, 1", $true, $false, $false, $false, $false, $true, 1, $false, "643÷5=128, 3", 2) | Out-Null
$d.Content.Find.Execute("313÷3=104, 1", $true, $false, $false, $false, $false, $true, 1, $false, "600÷6=100, 0", 2) | Out-Null
$d.Content.Find.Execute("589÷2=294, 1", $true, $false, $false, $false, $false, $true, 1, $false, "270÷3=90, 0", 2) | Out-Null
$d.Content.Find.Execute("789÷3=263, 0", $true, $false, $false, $false, $false, $true, 1, $false, "848÷9=94, 2", 2) | Out-Null
$d.Content.Find.Execute("115÷7=16, 3", $true, $false, $false, $false, $false, $true, 1, $false, "449÷2=224, 1", 2) | Out-Null
$d.Content.Find.Execute("100÷5=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "701÷6=116, 5", 2) | Out-Null
$d.Content.Find.Execute("893÷8=111, 5", $true, $false, $false, $false, $false, $true, 1, $false, "323÷4=80, 3", 2) | Out-Null
$d.Content.Find.Execute("349÷2=174, 1", $true, $false, $false, $false, $false, $true, 1, $false, "334÷2=167, 0", 2) | Out-Null
$d.Content.Find.Execute("697÷6=116, 1", $true, $false, $false, $false, $false, $true, 1, $false, "482÷5=96, 2", 2) | Out-Null
$d.Content.Find.Execute("826÷3=275, 1", $true, $false, $false, $false, $false, $true, 1, $false, "676÷2=338, 0", 2) | Out-Null
$d.Content.Find.Execute("456÷6=76, 0", $true, $false, $false, $false, $false, $true, 1, $false, "114÷9=12, 6", 2) | Out-Null
$d.Content.Find.Execute("341÷2=170, 1", $true, $false, $false, $false, $false, $true, 1, $false, "506÷5=101, 1", 2) | Out-Null
$d.Content.Find.Execute("414÷7=59, 1", $true, $false, $false, $false, $false, $true, 1, $false, "328÷3=109, 1", 2) | Out-Null
